# Fruta / hortaliza, semanal
# Re-order the weekly price records (rows 2-8) in columns D, K, L, M, N, O, P, Q, R, S, T
# according to the new mapping: new row N gets the values previously held by old row Map[N]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that participate in the record (everything that varies row to row)
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot current (pre-edit) values for rows 2-8 for each of those columns
$snapshot = @{}
foreach ($row in 2..8) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Mapping: new row -> source (old) row
$map = @{
    2 = 6
    3 = 4
    4 = 5
    5 = 2
    6 = 3
    7 = 8
    8 = 7
}

foreach ($newRow in 2..8) {
    $oldRow = $map[$newRow]
    $srcData = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $srcData[$col]
    }
}
